# BMB_Timesheet_Project.xlsx - add real timesheet entries for rows 4-16.
#
# Rows 3 stays as-is. Rows 4-6 get their Course-ID / Summary cells rewritten
# to the final (longer) course codes and new summary text; rows 7-16 (which
# were blank placeholder rows) get filled in with new entries. Finally the
# two content columns are resized and the selection is moved to B17 (the
# first still-empty row), matching where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pass 1: clear the cells whose text is being swapped for a brand-new
# string (not merely re-pointed) so the now-unused shared strings
# ("BUSA-245", "ECON-111", "PSYC") actually drop out of the shared string
# table instead of lingering as orphans. ---
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("F6").ClearContents()

# --- Pass 2: (re)write every data row, in sheet (row-major) order, so new
# shared strings are interned in the same order they are first encountered. ---

$ws.Range("D4").Value = "BUSA-245-03"
$ws.Range("F4").Value = "Class Disscussion"

$ws.Range("D5").Value = "ECON-111-01"
$ws.Range("F5").Value = "HW"

$ws.Range("D6").Value = "PSYC-103H-01"
$ws.Range("F6").Value = "Project"

$ws.Range("B7").Value = "BB"
$ws.Range("C7").Value = 45181
$ws.Range("C7").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D7").Value = "BUSA-245-03"
$ws.Range("E7").Value = 40
$ws.Range("F7").Value = "Take home quiz & disscussion post"

$ws.Range("B8").Value = "BB"
$ws.Range("C8").Value = 45181
$ws.Range("C8").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D8").Value = "DS-160-01"
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = "Finishing classwork"

$ws.Range("B9").Value = "BB"
$ws.Range("C9").Value = 45183
$ws.Range("C9").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D9").Value = "COMM-103-01"
$ws.Range("E9").Value = 20
$ws.Range("F9").Value = "Writing speech"

$ws.Range("B10").Value = "BB"
$ws.Range("C10").Value = 45186
$ws.Range("C10").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D10").Value = "BUSA-245-03"
$ws.Range("E10").Value = 30
$ws.Range("F10").Value = "Take home quiz & disscussion post"

$ws.Range("B11").Value = "BB"
$ws.Range("C11").Value = 45186
$ws.Range("C11").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D11").Value = "ECON-111-01"
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = "Cengage HW"

$ws.Range("B12").Value = "BB"
$ws.Range("C12").Value = 45189
$ws.Range("C12").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D12").Value = "DS-160-01"
$ws.Range("E12").Value = 60
$ws.Range("F12").Value = "HW"

$ws.Range("B13").Value = "BB"
$ws.Range("C13").Value = 45190
$ws.Range("C13").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D13").Value = "COMM-103-01"
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = "Practice Speech"

$ws.Range("B14").Value = "BB"
$ws.Range("C14").Value = 45190
$ws.Range("C14").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D14").Value = "ECON-111-01"
$ws.Range("E14").Value = 20
$ws.Range("F14").Value = "Cengage HW"

$ws.Range("B15").Value = "BB"
$ws.Range("C15").Value = 45193
$ws.Range("C15").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D15").Value = "PSYC-103H-01"
$ws.Range("E15").Value = 20
$ws.Range("F15").Value = "Project work"

$ws.Range("B16").Value = "BB"
$ws.Range("C16").Value = 45193
$ws.Range("C16").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D16").Value = "BUSA-245-03"
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = "Take home quiz"

# --- Column widths: Course-ID (D) best-fit-ish, Summary (F) widened. ---
$ws.Columns.Item(4).ColumnWidth = 11.830729166666666
$ws.Columns.Item(6).ColumnWidth = 18.166666666666668

# --- Move the selection to B17, the next empty row. ---
$ws.Range("B17").Select()
